$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 11797
$ws.Range("J17").Value = 11797
$ws.Range("L17").Value = 35391
$ws.Range("N17").Value = -35727
$ws.Range("H19").Value = 689.087
$ws.Range("I19").Value = 143.9
$ws.Range("K19").Value = 143.9
$ws.Range("M19").Value = 31.09999999999999
$ws.Range("H33").Value = 592.63635
$ws.Range("I33").Value = 671.8333
$ws.Range("J33").Value = 497.6
$ws.Range("K33").Value = 671.8333
$ws.Range("L33").Value = 497.6
$ws.Range("M33").Value = -442.8333
$ws.Range("N33").Value = -955.6
$ws.Range("H51").Value = 5563
$ws.Range("J51").Value = 5563
$ws.Range("L51").Value = 5563
$ws.Range("N51").Value = -6531
$ws.Range("H116").Value = 4917.1714
$ws.Range("I116").Value = 3788.1904
$ws.Range("J116").Value = 6610.643
$ws.Range("K116").Value = 3788.1904
$ws.Range("L116").Value = 6610.643
$ws.Range("M116").Value = -346.1904
$ws.Range("N116").Value = -13494.643
$ws.Range("H137").Value = 57888.25
$ws.Range("I137").Value = 70271
$ws.Range("J137").Value = 4229.6665
$ws.Range("K137").Value = 210813
$ws.Range("L137").Value = 12688.9995
$ws.Range("M137").Value = -208263
$ws.Range("N137").Value = -17788.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2357979.2
$ws.Range("I2").Value = 4714568
$ws.Range("J2").Value = 1390.6666
$ws.Range("K2").Value = 4714568
$ws.Range("L2").Value = 1390.6666
$ws.Range("M2").Value = -4714455
$ws.Range("N2").Value = -1616.6666
$ws.Range("H24").Value = 19677.5
$ws.Range("J24").Value = 19677.5
$ws.Range("L24").Value = 19677.5
$ws.Range("N24").Value = -20425.5
$ws.Range("H45").Value = 5921103.5
$ws.Range("I45").Value = 8099435.5
$ws.Range("J45").Value = 8488.857
$ws.Range("K45").Value = 8099435.5
$ws.Range("L45").Value = 8488.857
$ws.Range("M45").Value = -8099058.5
$ws.Range("N45").Value = -9242.857
$ws.Range("H55").Value = 39800
$ws.Range("I55").Value = 3000
$ws.Range("K55").Value = 3000
$ws.Range("M55").Value = -2685
$ws.Range("H88").Value = 1443.8125
$ws.Range("J88").Value = 1289.5555
$ws.Range("L88").Value = 1289.5555
$ws.Range("N88").Value = -2101.5555
$ws.Range("H91").Value = 1443.8125
$ws.Range("J91").Value = 1289.5555
$ws.Range("L91").Value = 1289.5555
$ws.Range("N91").Value = -4097.5555
$ws.Range("H97").Value = 1623576.4
$ws.Range("I97").Value = 3236135
$ws.Range("J97").Value = 11017.6
$ws.Range("K97").Value = 3236135
$ws.Range("L97").Value = 11017.6
$ws.Range("M97").Value = -3235639
$ws.Range("N97").Value = -12009.6
$ws.Range("H100").Value = 19677.5
$ws.Range("J100").Value = 19677.5
$ws.Range("L100").Value = 19677.5
$ws.Range("N100").Value = -21841.5
$ws.Range("H102").Value = 6415308
$ws.Range("I102").Value = 10421372
$ws.Range("J102").Value = 5606.2
$ws.Range("K102").Value = 10421372
$ws.Range("L102").Value = 5606.2
$ws.Range("M102").Value = -10419750
$ws.Range("N102").Value = -8850.200000000001
$ws.Range("H116").Value = 2357979.2
$ws.Range("I116").Value = 4714568
$ws.Range("J116").Value = 1390.6666
$ws.Range("K116").Value = 4714568
$ws.Range("L116").Value = 1390.6666
$ws.Range("M116").Value = -4712274
$ws.Range("N116").Value = -5978.6666
$ws.Range("H122").Value = 582969.8
$ws.Range("I122").Value = 3199
$ws.Range("K122").Value = 9597
$ws.Range("M122").Value = -7147

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2357979.2
$ws.Range("I3").Value = 4714568
$ws.Range("J3").Value = 1390.6666
$ws.Range("K3").Value = 4714568
$ws.Range("L3").Value = 1390.6666
$ws.Range("M3").Value = -4714454
$ws.Range("N3").Value = -1618.6666
$ws.Range("H22").Value = 1773.5264
$ws.Range("I22").Value = 1762.3125
$ws.Range("J22").Value = 1833.3334
$ws.Range("K22").Value = 1762.3125
$ws.Range("L22").Value = 1833.3334
$ws.Range("M22").Value = -1589.3125
$ws.Range("N22").Value = -2179.3334
$ws.Range("H86").Value = 5883624
$ws.Range("I86").Value = 11112233
$ws.Range("K86").Value = 11112233
$ws.Range("M86").Value = -11111110
$ws.Range("H89").Value = 5883624
$ws.Range("I89").Value = 11112233
$ws.Range("K89").Value = 55561165
$ws.Range("M89").Value = -55555549
$ws.Range("H99").Value = 6214964.5
$ws.Range("J99").Value = 3989.4443
$ws.Range("L99").Value = 3989.4443
$ws.Range("N99").Value = -6985.4443

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 3379.8
$ws.Range("I2").Value = 2974.75
$ws.Range("J2").Value = 5000
$ws.Range("K2").Value = 2974.75
$ws.Range("L2").Value = 5000
$ws.Range("M2").Value = -2861.75
$ws.Range("N2").Value = -5226
$ws.Range("H99").Value = 3123.5
$ws.Range("I99").Value = 2660.9375
$ws.Range("J99").Value = 4973.75
$ws.Range("K99").Value = 2660.9375
$ws.Range("L99").Value = 4973.75
$ws.Range("M99").Value = -1162.9375
$ws.Range("N99").Value = -7969.75
$ws.Range("H126").Value = 3123.5
$ws.Range("I126").Value = 2660.9375
$ws.Range("J126").Value = 4973.75
$ws.Range("K126").Value = 7982.8125
$ws.Range("L126").Value = 14921.25
$ws.Range("M126").Value = -5512.8125
$ws.Range("N126").Value = -19861.25
$ws.Range("H132").Value = 47175.56
$ws.Range("I132").Value = 30546.223
$ws.Range("K132").Value = 91638.66900000001
$ws.Range("M132").Value = -89108.66900000001
$ws.Range("H134").Value = 2944
$ws.Range("I134").Value = 2023.5
$ws.Range("K134").Value = 6070.5
$ws.Range("M134").Value = -3535.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 602.8823
$ws.Range("I25").Value = 556.6
$ws.Range("J25").Value = 950
$ws.Range("K25").Value = 1669.8
$ws.Range("L25").Value = 2850
$ws.Range("M25").Value = -1500.8
$ws.Range("N25").Value = -3188
$ws.Range("H30").Value = 602.8823
$ws.Range("I30").Value = 556.6
$ws.Range("J30").Value = 950
$ws.Range("K30").Value = 1669.8
$ws.Range("L30").Value = 2850
$ws.Range("M30").Value = -1567.8
$ws.Range("N30").Value = -3054
$ws.Range("H92").Value = 749.8461
$ws.Range("I92").Value = 282.33334
$ws.Range("J92").Value = 1150.5714
$ws.Range("K92").Value = 847.0000200000001
$ws.Range("L92").Value = 3451.7142
$ws.Range("M92").Value = 400.9999799999999
$ws.Range("N92").Value = -5947.7142
$ws.Range("H97").Value = 686.25
$ws.Range("J97").Value = 686.25
$ws.Range("L97").Value = 2058.75
$ws.Range("N97").Value = -3050.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 6176040.5
$ws.Range("I113").Value = 15153195
$ws.Range("J113").Value = 4246.75
$ws.Range("K113").Value = 15153195
$ws.Range("L113").Value = 4246.75
$ws.Range("M113").Value = -15151025
$ws.Range("N113").Value = -8586.75
$ws.Range("H122").Value = 288994.84
$ws.Range("I122").Value = 425210.16
$ws.Range("J122").Value = 2942.7
$ws.Range("K122").Value = 1275630.48
$ws.Range("L122").Value = 8828.099999999999
$ws.Range("M122").Value = -1273180.48
$ws.Range("N122").Value = -13728.1
$ws.Range("H132").Value = 3554.5
$ws.Range("I132").Value = 3336.5
$ws.Range("K132").Value = 10009.5
$ws.Range("M132").Value = -7479.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 7338.231
$ws.Range("I46").Value = 301
$ws.Range("K46").Value = 301
$ws.Range("M46").Value = -113
$ws.Range("H93").Value = 41682636
$ws.Range("I93").Value = 66669020
$ws.Range("J93").Value = 38666.332
$ws.Range("K93").Value = 66669020
$ws.Range("L93").Value = 38666.332
$ws.Range("M93").Value = -66667772
$ws.Range("N93").Value = -41162.332
$ws.Range("H132").Value = 6603.971
$ws.Range("I132").Value = 6786.3
$ws.Range("K132").Value = 20358.9
$ws.Range("M132").Value = -17828.9
$ws.Range("H136").Value = 28098.861
$ws.Range("I136").Value = 51037.715
$ws.Range("J136").Value = 6202.6816
$ws.Range("K136").Value = 153113.145
$ws.Range("L136").Value = 18608.0448
$ws.Range("M136").Value = -150563.145
$ws.Range("N136").Value = -23708.0448

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 43948.332
$ws.Range("J46").Value = 43948.332
$ws.Range("L46").Value = 43948.332
$ws.Range("N46").Value = -44410.332
$ws.Range("H81").Value = 9263526
$ws.Range("I81").Value = 16668617
$ws.Range("K81").Value = 33337234
$ws.Range("M81").Value = -33336173
$ws.Range("H84").Value = 9263526
$ws.Range("I84").Value = 16668617
$ws.Range("K84").Value = 166686170
$ws.Range("M84").Value = -166680866
$ws.Range("H93").Value = 49999.5
$ws.Range("J93").Value = 49999.5
$ws.Range("L93").Value = 49999.5
$ws.Range("N93").Value = -54991.5
$ws.Range("H134").Value = 43948.332
$ws.Range("J134").Value = 43948.332
$ws.Range("L134").Value = 131844.996
$ws.Range("N134").Value = -136914.996
